# "updated with quiz 5"
# Populates the QUIZ FIVE (column J) scores on the "Class Quizes" sheet,
# fills a couple of late-added LAB TWO (column H) / LAB ONE (column G) scores,
# and appends a trailing row with a lone "\" marker in column L.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Class Quizes")

# --- Quiz 5 (column J) scores, plus a few missed Lab scores (columns G/H) ---
$ws.Range("J2").Value = 4
$ws.Range("J3").Value = 4
$ws.Range("J5").Value = 10
$ws.Range("J6").Value = 1
$ws.Range("J7").Value = 5
$ws.Range("J8").Value = 9
$ws.Range("H9").Value = 4
$ws.Range("J9").Value = 7
$ws.Range("J10").Value = 5
$ws.Range("J11").Value = 5
$ws.Range("J12").Value = 6
$ws.Range("J13").Value = 5
$ws.Range("J15").Value = 4
$ws.Range("J16").Value = 4
$ws.Range("J17").Value = 10
$ws.Range("J18").Value = 4
$ws.Range("J20").Value = 3
$ws.Range("J21").Value = 10
$ws.Range("H22").Value = 4
$ws.Range("J22").Value = 6
$ws.Range("J23").Value = 4
$ws.Range("G24").Value = 3
$ws.Range("H24").Value = 4
$ws.Range("J24").Value = 4
$ws.Range("J25").Value = 4
$ws.Range("J26").Value = 9
$ws.Range("J28").Value = 9
$ws.Range("J29").Value = 4
$ws.Range("J30").Value = 4
$ws.Range("H31").Value = 4
$ws.Range("J32").Value = 6
$ws.Range("J33").Value = 4
$ws.Range("J35").Value = 5
$ws.Range("J36").Value = 4
$ws.Range("J37").Value = 6
$ws.Range("J38").Value = 4
$ws.Range("J39").Value = 4
$ws.Range("J40").Value = 7
$ws.Range("J42").Value = 5
$ws.Range("J43").Value = 4
$ws.Range("J44").Value = 3
$ws.Range("J46").Value = 5
$ws.Range("J47").Value = 7
$ws.Range("J52").Value = 5

# --- New trailing row with a lone backslash marker ---
$ws.Range("L54").Value = "\"

# --- Restore the view position/selection as left by the editor ---
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 11
$excel.ActiveWindow.ScrollRow = 32
[void]$ws.Range("M38").Select()
